$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4216.6665
$ws.Range("J29").Value = 4216.6665
$ws.Range("L29").Value = 12649.9995
$ws.Range("N29").Value = -13211.9995
$ws.Range("H33").Value = 559.7143
$ws.Range("I33").Value = 425
$ws.Range("K33").Value = 425
$ws.Range("M33").Value = -196
$ws.Range("H38").Value = 816.3077
$ws.Range("I38").Value = 296.38095
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 889.14285
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -517.14285
$ws.Range("N38").Value = -9744
$ws.Range("H58").Value = 40287.69
$ws.Range("I58").Value = 942
$ws.Range("J58").Value = 64878.75
$ws.Range("K58").Value = 2826
$ws.Range("L58").Value = 194636.25
$ws.Range("M58").Value = -2676
$ws.Range("N58").Value = -194936.25
$ws.Range("H69").Value = 5310.75
$ws.Range("J69").Value = 5910
$ws.Range("L69").Value = 17730
$ws.Range("N69").Value = -19478
$ws.Range("H72").Value = 5310.75
$ws.Range("J72").Value = 5910
$ws.Range("L72").Value = 53190
$ws.Range("N72").Value = -61926
$ws.Range("H82").Value = 6355.4
$ws.Range("J82").Value = 7819.25
$ws.Range("L82").Value = 23457.75
$ws.Range("N82").Value = -24269.75
$ws.Range("H85").Value = 6355.4
$ws.Range("J85").Value = 7819.25
$ws.Range("L85").Value = 23457.75
$ws.Range("N85").Value = -26265.75
$ws.Range("H87").Value = 25747
$ws.Range("J87").Value = 25747
$ws.Range("L87").Value = 25747
$ws.Range("N87").Value = -28243
$ws.Range("H90").Value = 25747
$ws.Range("J90").Value = 25747
$ws.Range("L90").Value = 77241
$ws.Range("N90").Value = -89721
$ws.Range("H96").Value = 706
$ws.Range("I96").Value = 778.6667
$ws.Range("J96").Value = 618.8
$ws.Range("K96").Value = 2336.0001
$ws.Range("L96").Value = 1856.4
$ws.Range("M96").Value = -963.0001000000002
$ws.Range("N96").Value = -4602.4
$ws.Range("H129").Value = 987.1395
$ws.Range("I129").Value = 941.8182
$ws.Range("J129").Value = 1002.71875
$ws.Range("K129").Value = 2825.4546
$ws.Range("L129").Value = 3008.15625
$ws.Range("M129").Value = 2174.5454
$ws.Range("N129").Value = -13008.15625
$ws.Range("H131").Value = 3459.7273
$ws.Range("I131").Value = 3492.6
$ws.Range("J131").Value = 3432.3333
$ws.Range("K131").Value = 10477.8
$ws.Range("L131").Value = 10296.9999
$ws.Range("M131").Value = -5437.799999999999
$ws.Range("N131").Value = -20376.9999
$ws.Range("H135").Value = 592.1429000000001
$ws.Range("I135").Value = 460.4
$ws.Range("J135").Value = 1690
$ws.Range("K135").Value = 4143.599999999999
$ws.Range("L135").Value = 15210
$ws.Range("M135").Value = -1608.599999999999
$ws.Range("N135").Value = -20280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 21344
$ws.Range("J37").Value = 23250.285
$ws.Range("L37").Value = 23250.285
$ws.Range("N37").Value = -23796.285
$ws.Range("H44").Value = 21308.8
$ws.Range("J44").Value = 21308.8
$ws.Range("L44").Value = 21308.8
$ws.Range("N44").Value = -22284.8
$ws.Range("H55").Value = 25053
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 25053
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = 25053
$ws.Range("N55").Value = -25683
$ws.Range("H80").Value = 27124
$ws.Range("J80").Value = 27124
$ws.Range("L80").Value = 27124
$ws.Range("N80").Value = -29120
$ws.Range("H83").Value = 27124
$ws.Range("J83").Value = 27124
$ws.Range("L83").Value = 81372
$ws.Range("N83").Value = -91356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 26740.666
$ws.Range("J35").Value = 26740.666
$ws.Range("L35").Value = 26740.666
$ws.Range("N35").Value = -27360.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 6898.4546
$ws.Range("I41").Value = 2726
$ws.Range("J41").Value = 11905.4
$ws.Range("K41").Value = 2726
$ws.Range("L41").Value = 11905.4
$ws.Range("M41").Value = -2298
$ws.Range("N41").Value = -12761.4
$ws.Range("H50").Value = 9225.143
$ws.Range("J50").Value = 9225.143
$ws.Range("L50").Value = 9225.143
$ws.Range("N50").Value = -10475.143
$ws.Range("H51").Value = 18056.857
$ws.Range("I51").Value = 9000
$ws.Range("J51").Value = 19566.334
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 19566.334
$ws.Range("M51").Value = -8264
$ws.Range("N51").Value = -21038.334
$ws.Range("H61").Value = 18056.857
$ws.Range("I61").Value = 9000
$ws.Range("J61").Value = 19566.334
$ws.Range("K61").Value = 9000
$ws.Range("L61").Value = 19566.334
$ws.Range("M61").Value = -8652
$ws.Range("N61").Value = -20262.334
$ws.Range("H74").Value = 15773.846
$ws.Range("J74").Value = 15773.846
$ws.Range("L74").Value = 15773.846
$ws.Range("N74").Value = -17521.846
$ws.Range("H77").Value = 15773.846
$ws.Range("J77").Value = 15773.846
$ws.Range("L77").Value = 47321.538
$ws.Range("N77").Value = -56057.538
$ws.Range("H132").Value = 4401.5713
$ws.Range("I132").Value = 4703.4287
$ws.Range("J132").Value = 4099.7144
$ws.Range("K132").Value = 14110.2861
$ws.Range("L132").Value = 12299.1432
$ws.Range("M132").Value = -11580.2861
$ws.Range("N132").Value = -17359.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6949.5884
$ws.Range("J34").Value = 11691.3
$ws.Range("L34").Value = 35073.89999999999
$ws.Range("N34").Value = -35241.89999999999
$ws.Range("H39").Value = 2005.6786
$ws.Range("I39").Value = 750
$ws.Range("J39").Value = 2102.2693
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 6306.8079
$ws.Range("M39").Value = -1956
$ws.Range("N39").Value = -6894.8079
$ws.Range("H55").Value = 2017.7142
$ws.Range("I55").Value = 514
$ws.Range("J55").Value = 3521.4285
$ws.Range("K55").Value = 1542
$ws.Range("L55").Value = 10564.2855
$ws.Range("M55").Value = -1365
$ws.Range("N55").Value = -10918.2855
$ws.Range("H86").Value = 497.34784
$ws.Range("I86").Value = 79
$ws.Range("K86").Value = 237
$ws.Range("M86").Value = 949
$ws.Range("H89").Value = 497.34784
$ws.Range("I89").Value = 79
$ws.Range("K89").Value = 711
$ws.Range("M89").Value = 5217
$ws.Range("H131").Value = 1256.8903
$ws.Range("I131").Value = 4672
$ws.Range("J131").Value = 1035.1299
$ws.Range("K131").Value = 14016
$ws.Range("L131").Value = 3105.3897
$ws.Range("M131").Value = -8976
$ws.Range("N131").Value = -13185.3897
$ws.Range("H137").Value = 2172.139
$ws.Range("I137").Value = 1576.5769
$ws.Range("J137").Value = 3720.6
$ws.Range("K137").Value = 4729.7307
$ws.Range("L137").Value = 11161.8
$ws.Range("M137").Value = 370.2692999999999
$ws.Range("N137").Value = -21361.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 27145.56
$ws.Range("I102").Value = 1957.4138
$ws.Range("J102").Value = 88016.914
$ws.Range("K102").Value = 1957.4138
$ws.Range("L102").Value = 88016.914
$ws.Range("M102").Value = -335.4138
$ws.Range("N102").Value = -91260.914
$ws.Range("H132").Value = 4195.875
$ws.Range("I132").Value = 4605.625
$ws.Range("K132").Value = 13816.875
$ws.Range("M132").Value = -11286.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 11699.111
$ws.Range("I35").Value = 1058.4
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 1058.4
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -722.4000000000001
$ws.Range("N35").Value = -25672

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14007.75
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 17010.334
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 17010.334
$ws.Range("M54").Value = -4480
$ws.Range("N54").Value = -18050.334
$ws.Range("H81").Value = 1335.7693
$ws.Range("I81").Value = 1139.5
$ws.Range("J81").Value = 1990
$ws.Range("K81").Value = 2279
$ws.Range("L81").Value = 3980
$ws.Range("M81").Value = -1218
$ws.Range("N81").Value = -6102
$ws.Range("H84").Value = 1335.7693
$ws.Range("I84").Value = 1139.5
$ws.Range("J84").Value = 1990
$ws.Range("K84").Value = 11395
$ws.Range("L84").Value = 19900
$ws.Range("M84").Value = -6091
$ws.Range("N84").Value = -30508

